$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Helper: find the paragraph whose text starts with $marker and replace its
# whole XML with $newXml (an OOXML <w:p>...</w:p> fragment, without the
# wordprocessingml namespace declaration - that gets spliced in).
function Replace-ParagraphXml {
    param(
        [string]$Marker,
        [string]$InnerXml
    )
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($Marker)) {
            $full = $InnerXml -replace '<w:p ', ("<w:p $wns ")
            if ($full -eq $InnerXml) {
                $full = $InnerXml -replace '<w:p>', ("<w:p $wns>")
            }
            $p.Range.InsertXML($full)
            return
        }
    }
    throw "Paragraph starting with '$Marker' not found"
}

# 1) "Fri. 11/21" heading paragraph: drop <w:rFonts w:hint="eastAsia"/>
#    from the paragraph-mark's rPr (w:pPr/w:rPr), keep everything else.
Replace-ParagraphXml "Fri. 11/21" (
    "<w:p w:rsidR='001A7EC6' w:rsidRPr='001A7EC6' w:rsidRDefault='001A7EC6'>" +
    "<w:pPr><w:rPr><w:u w:val='single'/></w:rPr></w:pPr>" +
    "<w:r w:rsidRPr='001A7EC6'><w:rPr><w:rFonts w:hint='eastAsia'/><w:u w:val='single'/></w:rPr><w:t>Fri. 11/21</w:t></w:r>" +
    "</w:p>"
)

# 2) "Sat. 11/22" heading paragraph: same rFonts removal.
Replace-ParagraphXml "Sat. 11/22" (
    "<w:p w:rsidR='001A7EC6' w:rsidRDefault='001A7EC6' w:rsidP='001A7EC6'>" +
    "<w:pPr><w:rPr><w:u w:val='single'/></w:rPr></w:pPr>" +
    "<w:r w:rsidRPr='001A7EC6'><w:rPr><w:rFonts w:hint='eastAsia'/><w:u w:val='single'/></w:rPr><w:t>Sat. 11/22</w:t></w:r>" +
    "</w:p>"
)

# 3) "Comment each Code" paragraph: drop the pPr (only held rFonts) and
#    rewrite the single run as four runs: "DONE-", "Comment ",
#    "L1Cache and ", and a spell-checked "MainMemory".
Replace-ParagraphXml "Comment each Code" (
    "<w:p w:rsidR='001A7EC6' w:rsidRDefault='001A7EC6' w:rsidP='001A7EC6'>" +
    "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>DONE-</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t xml:space='preserve'>Comment </w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t xml:space='preserve'>L1Cache and </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>MainMemory</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
)

# 4) "Read cache line from DRAM" list paragraph: drop the trailing
#    <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> that sat inside pPr
#    after numPr (paragraph-mark run properties), keep numPr/pStyle.
Replace-ParagraphXml "Read cache line from DRAM" (
    "<w:p w:rsidR='001A7EC6' w:rsidRDefault='001A7EC6' w:rsidP='001A7EC6'>" +
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/></w:rPr><w:t>Read cache line from DRAM</w:t></w:r>" +
    "</w:p>"
)
